# Refresh the cryptos price list (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are refreshed for most rows; a
# couple of rows also got their Coin name/Link swapped back into the
# "correct" rank order (Uniswap/InternetComputer(DFINITY) at 19/20 and
# Celestia/LidoDAOToken at 34/35).
#
# Some "Price" strings (e.g. "315.53") look like plain numbers to Excel,
# which would silently convert them to numeric values on assignment. The
# source data stores every Price/Volume cell as text, so for any D-column
# value that parses as a number we briefly force a text NumberFormat,
# assign the value, then restore the cell's style to "Normal" so no stray
# formatting is left behind (values with extra "." like "42.839.80" are
# never auto-converted, so they're set directly).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.839.80'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '2.532.15'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.577'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.62%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.83%  '

$ws.Range("E11").Value = '  -0.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.110'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.86%  '

$ws.Range("D14").Value = '2.919.42'
$ws.Range("E14").Value = '  +0.34%  '

$ws.Range("D15").Value = '2.603.78'
$ws.Range("E15").Value = '  +4.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.852'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.53%  '

$ws.Range("D18").Value = '42.885.56'
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.00%  '

$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.65%  '

$ws.Range("E21").Value = '  -1.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.58%  '

$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.64'
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("E28").Value = '  +2.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.17%  '

$ws.Range("E33").Value = '  +3.43%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.13'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0782'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.98%  '

$ws.Range("E38").Value = '  -1.77%  '

$ws.Range("E39").Value = '  -1.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.35%  '

$ws.Range("E41").Value = '  +13.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("E45").Value = '  -2.81%  '

$ws.Range("D46").Value = '2.029.56'
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.41'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.33%  '

$ws.Range("D51").Value = '2.773.41'
$ws.Range("E51").Value = '  +0.35%  '
